# feat: add 2022-Q1 data
#
# 1) Duplicate the "2021-Q4" sheet (which already carries the per-fund table
#    layout/styles) to create the new "2022-Q1" sheet right after it, then
#    overwrite its two data rows with the 2022-Q1 fund holdings.
# 2) Insert a new top data row into the "总计" (totals) sheet for 2022-Q1,
#    pushing the existing 2021-Q4 totals row down to row 3.

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, [string]$text) {
    # Force text storage (not numeric auto-detection) the way Excel does
    # for a user-typed quote-prefixed entry, then drop back to the
    # worksheet's default style so no stray number-format/quote-prefix
    # style sticks around on the cell.
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# ------------------------------------------------------------------
# 1) "2022-Q1" sheet: clone "2021-Q4" immediately after itself
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)

$q1 = $wb.Worksheets.Item(2)
$q1.Name = "2022-Q1"

# Row 2: 富国中证价值ETF
Set-TextValue $q1.Range("B2") "512040"
Set-TextValue $q1.Range("C2") "富国中证价值ETF"
Set-TextValue $q1.Range("D2") "3.44"
Set-TextValue $q1.Range("E2") "99.55"
Set-TextValue $q1.Range("F2") "1.51"
Set-TextValue $q1.Range("G2") "0.0519"
$q1.Range("H2").Value = 3

# Row 3: 银河定投宝中证腾讯济安价值100A股指数
Set-TextValue $q1.Range("B3") "519677"
Set-TextValue $q1.Range("C3") "银河定投宝中证腾讯济安价值100A股指数"
Set-TextValue $q1.Range("D3") "2.74"
Set-TextValue $q1.Range("E3") "91.56"
Set-TextValue $q1.Range("F3") "1.35"
Set-TextValue $q1.Range("G3") "0.0370"
$q1.Range("H3").Value = 3

# ------------------------------------------------------------------
# 2) "总计" sheet: add the 2022-Q1 totals row above the 2021-Q4 one
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Copy column-A's styling down to the new row 3 before filling it in.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$total.Range("A3").Value = 1
Set-TextValue $total.Range("B3") "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.06

$total.Range("A2").Value = 0
Set-TextValue $total.Range("B2") "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.09

# Keep the originally active "2021-Q4" tab selected/active.
$q4.Activate()

Write-Host "2022-Q1 data added"
